$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Value 1"/"Value 2"/"Value 3" header cells to "Data1"/"Data2"/"Data3"
$ws.Range("B1").Value = "Data1"
$ws.Range("C1").Value = "Data2"
$ws.Range("D1").Value = "Data3"

# Move the current selection/active cell to D1
$ws.Range("D1").Select()
